$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Insert a new row right after row 13 (pushes the old B14/B15 bullets down to B15/B16)
$ws.Rows(14).Insert()

# Reword the existing bullet in B13, and give the freshly inserted B14 the new continuation sentence
$ws.Range("B13").Value = "* Do not use a formula in a cell that may have its position changed after the merge (for example under a TBS block). Otherwise Excel will raise an error message."
$ws.Range("B14").Value = "    This is because the location of formulas are saved a second time in another sub-file for the order of evaluation."

# B15/B16 already hold the old B14/B15 bullets verbatim thanks to the row insert above.
# Row 19 (old blank row between "Example #1" and the "First Name" header) is already free.

# New "Total:" label + SUM formula on row 19 (columns D/E)
$ws.Range("D19").Value = "Total:"
$ws.Range("D19").HorizontalAlignment = -4152
$ws.Range("E19").Formula = "=SUM(E21:E2000)"
$ws.Range("E19").Font.Bold = $true
$ws.Range("E19").NumberFormat = "#,##0.0"

# New "Score" column header + example field
$ws.Range("E20").Value = "Score"
$ws.Range("E20").Interior.Pattern = -4124
$ws.Range("E20").Interior.PatternColorIndex = -4105
$ws.Range("E20").Interior.ThemeColor = 0
$ws.Range("E20").Interior.TintAndShade = -0.14999847407452621
$ws.Range("E20").Borders.LineStyle = 1

$ws.Range("E21").Value = "[a.score;ope=xlsxNum]"
$ws.Range("E21").NumberFormat = "#,##0.0"
$ws.Range("E21").HorizontalAlignment = -4152
$ws.Range("E21").Borders.LineStyle = 1

$ws.Range("E20").Select()
